$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to Text format so numeric-looking strings
# (e.g. "43.280.36", "0.870", "1.00") are preserved exactly as text
# instead of being parsed/rounded as numbers.

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "43.280.36"
$ws.Cells.Item(2, 5).Value = "  -4.58%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.240.51"
$ws.Cells.Item(3, 5).Value = "  -5.78%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.24%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "321.34"
$ws.Cells.Item(5, 5).Value = "  +1.27%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "101.62"
$ws.Cells.Item(6, 5).Value = "  -6.53%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.587"
$ws.Cells.Item(7, 5).Value = "  -8.17%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.15%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.565"
$ws.Cells.Item(9, 5).Value = "  -8.42%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "37.13"
$ws.Cells.Item(10, 5).Value = "  -9.23%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "54.59"
$ws.Cells.Item(11, 5).Value = "  -2.56%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -9.84%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "7.74"
$ws.Cells.Item(13, 5).Value = "  -9.30%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.109"
$ws.Cells.Item(14, 5).Value = "  -0.82%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "Polygon"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.870"
$ws.Cells.Item(15, 5).Value = "  -11.82%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.580.68"
$ws.Cells.Item(16, 5).Value = "  -5.82%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.44"
$ws.Cells.Item(17, 5).Value = "  -7.31%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.250.67"
$ws.Cells.Item(18, 5).Value = "  -5.18%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "43.239.41"
$ws.Cells.Item(19, 5).Value = "  -4.52%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.44"
$ws.Cells.Item(20, 5).Value = "  -7.88%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.0₃0969"
$ws.Cells.Item(21, 5).Value = "  -8.94%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.55"
$ws.Cells.Item(22, 5).Value = "  -10.50%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "65.59"
$ws.Cells.Item(23, 5).Value = "  -10.56%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -13.10%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "238.14"
$ws.Cells.Item(25, 5).Value = "  -8.86%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.17"
$ws.Cells.Item(26, 5).Value = "  -7.62%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.00"
$ws.Cells.Item(27, 5).Value = "  -0.06%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "4.06"
$ws.Cells.Item(28, 5).Value = "  +1.65%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "Cosmos"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "10.07"
$ws.Cells.Item(29, 5).Value = "  -10.43%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Toncoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.18"
$ws.Cells.Item(30, 5).Value = "  -4.73%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -16.11%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "36.17"
$ws.Cells.Item(32, 5).Value = "  -2.96%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0883"
$ws.Cells.Item(33, 5).Value = "  -8.07%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "EthereumClassic"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "20.48"
$ws.Cells.Item(34, 5).Value = "  -8.47%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "153.91"
$ws.Cells.Item(35, 5).Value = "  -8.09%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -4.84%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.25"
$ws.Cells.Item(37, 5).Value = "  +9.27%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +1.92%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -7.81%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "4.47"
$ws.Cells.Item(40, 5).Value = "  -5.24%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -10.75%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.72"
$ws.Cells.Item(42, 5).Value = "  -8.13%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0326"
$ws.Cells.Item(43, 5).Value = "  -8.12%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "13.49"
$ws.Cells.Item(44, 5).Value = "  +3.42%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -0.09%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.775.31"
$ws.Cells.Item(46, 5).Value = "  -4.02%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "87.55"
$ws.Cells.Item(47, 5).Value = "  -10.59%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -10.24%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "ordi"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "76.73"
$ws.Cells.Item(49, 5).Value = "  -8.20%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "THORChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.34"
$ws.Cells.Item(50, 5).Value = "  -11.06%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "59.32"
$ws.Cells.Item(51, 5).Value = "  -15.81%  "
